# This script applies the diff changes to the Denmark 2nd-division 2023-2024
# betting-odds workbook:
#   1. Rows 5 & 6, 14 & 15, 26 & 27 have their match data (columns F:V) swapped.
#   2. Rows 28/29/30 are cyclically rotated (28<-29, 29<-30, 30<-28).
#   3. Rows 31/32/33 are cyclically rotated (31<-32, 32<-33, 33<-31).
#   4. A brand-new match row (row 43 / Indice 42, Roskilde vs Hellerup) is
#      appended after the existing last row (42), extending the used range
#      from A1:V42 to A1:V43.
# Columns A-E (Indice, pais, torneio, temporada, data_partida) are untouched
# for the swapped/rotated rows; only F:V (home..url_partida) move.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 <- old row 6 data
$ws.Cells.Item(5,6).Value = "Brabrand"
$ws.Cells.Item(5,7).Value = 1
$ws.Cells.Item(5,8).Value = "Thisted FC"
$ws.Cells.Item(5,9).Value = 2
$ws.Cells.Item(5,10).Value = 1.95
$ws.Cells.Item(5,11).Value = "04/08/2023 02:12"
$ws.Cells.Item(5,12).Value = 2.47
$ws.Cells.Item(5,13).Value = "05/08/2023 12:25"
$ws.Cells.Item(5,14).Value = 3.36
$ws.Cells.Item(5,15).Value = "04/08/2023 02:12"
$ws.Cells.Item(5,16).Value = 3.4
$ws.Cells.Item(5,17).Value = "05/08/2023 12:04"
$ws.Cells.Item(5,18).Value = 3.26
$ws.Cells.Item(5,19).Value = "04/08/2023 02:12"
$ws.Cells.Item(5,20).Value = 2.69
$ws.Cells.Item(5,21).Value = "05/08/2023 12:25"
$ws.Cells.Item(5,22).Value = "https://www.betexplorer.com/football/denmark/2nd-division/brabrand-thisted-fc/bg3D6Tus/"
# Row 6 <- old row 5 data
$ws.Cells.Item(6,6).Value = "Roskilde"
$ws.Cells.Item(6,7).Value = 1
$ws.Cells.Item(6,8).Value = "FA 2000"
$ws.Cells.Item(6,9).Value = 0
$ws.Cells.Item(6,10).Value = 1.58
$ws.Cells.Item(6,11).Value = "04/08/2023 22:16"
$ws.Cells.Item(6,12).Value = 1.61
$ws.Cells.Item(6,13).Value = "05/08/2023 13:01"
$ws.Cells.Item(6,14).Value = 4.16
$ws.Cells.Item(6,15).Value = "04/08/2023 22:16"
$ws.Cells.Item(6,16).Value = 4.53
$ws.Cells.Item(6,17).Value = "05/08/2023 13:01"
$ws.Cells.Item(6,18).Value = 4.83
$ws.Cells.Item(6,19).Value = "04/08/2023 22:16"
$ws.Cells.Item(6,20).Value = 4.34
$ws.Cells.Item(6,21).Value = "05/08/2023 13:01"
$ws.Cells.Item(6,22).Value = "https://www.betexplorer.com/football/denmark/2nd-division/roskilde-frederiksberg-alliancen-2000/fLdM47Ag/"
# Row 14 <- old row 15 data
$ws.Cells.Item(14,6).Value = "AB Copenhagen"
$ws.Cells.Item(14,7).Value = 5
$ws.Cells.Item(14,8).Value = "Hellerup"
$ws.Cells.Item(14,9).Value = 4
$ws.Cells.Item(14,10).Value = 2.01
$ws.Cells.Item(14,11).Value = "17/08/2023 07:12"
$ws.Cells.Item(14,12).Value = 1.97
$ws.Cells.Item(14,13).Value = "18/08/2023 18:50"
$ws.Cells.Item(14,14).Value = 3.53
$ws.Cells.Item(14,15).Value = "17/08/2023 07:12"
$ws.Cells.Item(14,16).Value = 3.65
$ws.Cells.Item(14,17).Value = "18/08/2023 18:50"
$ws.Cells.Item(14,18).Value = 3
$ws.Cells.Item(14,19).Value = "17/08/2023 07:12"
$ws.Cells.Item(14,20).Value = 3.46
$ws.Cells.Item(14,21).Value = "18/08/2023 18:50"
$ws.Cells.Item(14,22).Value = "https://www.betexplorer.com/football/denmark/2nd-division/ab-copenhagen-hellerup/86V2Hpfa/"
# Row 15 <- old row 14 data
$ws.Cells.Item(15,6).Value = "Nykobing"
$ws.Cells.Item(15,7).Value = 1
$ws.Cells.Item(15,8).Value = "Roskilde"
$ws.Cells.Item(15,9).Value = 2
$ws.Cells.Item(15,10).Value = 2.04
$ws.Cells.Item(15,11).Value = "17/08/2023 21:49"
$ws.Cells.Item(15,12).Value = 2.38
$ws.Cells.Item(15,13).Value = "18/08/2023 17:37"
$ws.Cells.Item(15,14).Value = 3.5
$ws.Cells.Item(15,15).Value = "17/08/2023 21:49"
$ws.Cells.Item(15,16).Value = 3.81
$ws.Cells.Item(15,17).Value = "18/08/2023 17:34"
$ws.Cells.Item(15,18).Value = 3.31
$ws.Cells.Item(15,19).Value = "17/08/2023 21:49"
$ws.Cells.Item(15,20).Value = 2.57
$ws.Cells.Item(15,21).Value = "18/08/2023 17:37"
$ws.Cells.Item(15,22).Value = "https://www.betexplorer.com/football/denmark/2nd-division/nykobing-roskilde/KtZ6G495/"
# Row 26 <- old row 27 data
$ws.Cells.Item(26,6).Value = "Brabrand"
$ws.Cells.Item(26,7).Value = 1
$ws.Cells.Item(26,8).Value = "Hellerup"
$ws.Cells.Item(26,9).Value = 2
$ws.Cells.Item(26,10).Value = 2.51
$ws.Cells.Item(26,11).Value = "01/09/2023 02:12"
$ws.Cells.Item(26,12).Value = 2.37
$ws.Cells.Item(26,13).Value = "01/09/2023 23:41"
$ws.Cells.Item(26,14).Value = 3.34
$ws.Cells.Item(26,15).Value = "01/09/2023 02:12"
$ws.Cells.Item(26,16).Value = 3.56
$ws.Cells.Item(26,17).Value = "02/09/2023 12:03"
$ws.Cells.Item(26,18).Value = 2.45
$ws.Cells.Item(26,19).Value = "01/09/2023 02:12"
$ws.Cells.Item(26,20).Value = 2.68
$ws.Cells.Item(26,21).Value = "01/09/2023 23:41"
$ws.Cells.Item(26,22).Value = "https://www.betexplorer.com/football/denmark/2nd-division/brabrand-hellerup/EP9oS18H/"
# Row 27 <- old row 26 data
$ws.Cells.Item(27,6).Value = "Skive"
$ws.Cells.Item(27,7).Value = 2
$ws.Cells.Item(27,8).Value = "FA 2000"
$ws.Cells.Item(27,9).Value = 2
$ws.Cells.Item(27,10).Value = 2.44
$ws.Cells.Item(27,11).Value = "02/09/2023 10:43"
$ws.Cells.Item(27,12).Value = 2.58
$ws.Cells.Item(27,13).Value = "02/09/2023 13:48"
$ws.Cells.Item(27,14).Value = 3.61
$ws.Cells.Item(27,15).Value = "02/09/2023 10:43"
$ws.Cells.Item(27,16).Value = 3.43
$ws.Cells.Item(27,17).Value = "02/09/2023 13:48"
$ws.Cells.Item(27,18).Value = 2.56
$ws.Cells.Item(27,19).Value = "02/09/2023 10:43"
$ws.Cells.Item(27,20).Value = 2.55
$ws.Cells.Item(27,21).Value = "02/09/2023 13:48"
$ws.Cells.Item(27,22).Value = "https://www.betexplorer.com/football/denmark/2nd-division/skive-frederiksberg-alliancen-2000/K2AsTsgB/"
# Row 28 <- old row 29 data
$ws.Cells.Item(28,6).Value = "Nykobing"
$ws.Cells.Item(28,7).Value = 0
$ws.Cells.Item(28,8).Value = "Middelfart"
$ws.Cells.Item(28,9).Value = 0
$ws.Cells.Item(28,10).Value = 2.81
$ws.Cells.Item(28,11).Value = "02/09/2023 14:13"
$ws.Cells.Item(28,12).Value = 2.48
$ws.Cells.Item(28,13).Value = "02/09/2023 14:44"
$ws.Cells.Item(28,14).Value = 3.81
$ws.Cells.Item(28,15).Value = "02/09/2023 14:13"
$ws.Cells.Item(28,16).Value = 3.64
$ws.Cells.Item(28,17).Value = "02/09/2023 14:42"
$ws.Cells.Item(28,18).Value = 2.21
$ws.Cells.Item(28,19).Value = "02/09/2023 14:13"
$ws.Cells.Item(28,20).Value = 2.54
$ws.Cells.Item(28,21).Value = "02/09/2023 14:44"
$ws.Cells.Item(28,22).Value = "https://www.betexplorer.com/football/denmark/2nd-division/nykobing-middelfart/Ai3fQuwU/"
# Row 29 <- old row 30 data
$ws.Cells.Item(29,6).Value = "AB Copenhagen"
$ws.Cells.Item(29,7).Value = 0
$ws.Cells.Item(29,8).Value = "Esbjerg"
$ws.Cells.Item(29,9).Value = 3
$ws.Cells.Item(29,10).Value = 3.69
$ws.Cells.Item(29,11).Value = "01/09/2023 03:12"
$ws.Cells.Item(29,12).Value = 4.04
$ws.Cells.Item(29,13).Value = "02/09/2023 14:59"
$ws.Cells.Item(29,14).Value = 3.86
$ws.Cells.Item(29,15).Value = "01/09/2023 03:12"
$ws.Cells.Item(29,16).Value = 4.18
$ws.Cells.Item(29,17).Value = "02/09/2023 14:59"
$ws.Cells.Item(29,18).Value = 1.71
$ws.Cells.Item(29,19).Value = "01/09/2023 03:12"
$ws.Cells.Item(29,20).Value = 1.71
$ws.Cells.Item(29,21).Value = "02/09/2023 14:59"
$ws.Cells.Item(29,22).Value = "https://www.betexplorer.com/football/denmark/2nd-division/ab-copenhagen-esbjerg/SSDkRLNN/"
# Row 30 <- old row 28 data
$ws.Cells.Item(30,6).Value = "Thisted FC"
$ws.Cells.Item(30,7).Value = 1
$ws.Cells.Item(30,8).Value = "Aarhus Fremad"
$ws.Cells.Item(30,9).Value = 3
$ws.Cells.Item(30,10).Value = 4.21
$ws.Cells.Item(30,11).Value = "01/09/2023 03:12"
$ws.Cells.Item(30,12).Value = 4.52
$ws.Cells.Item(30,13).Value = "02/09/2023 14:40"
$ws.Cells.Item(30,14).Value = 3.95
$ws.Cells.Item(30,15).Value = "01/09/2023 03:12"
$ws.Cells.Item(30,16).Value = 3.95
$ws.Cells.Item(30,17).Value = "02/09/2023 14:51"
$ws.Cells.Item(30,18).Value = 1.6
$ws.Cells.Item(30,19).Value = "01/09/2023 03:12"
$ws.Cells.Item(30,20).Value = 1.68
$ws.Cells.Item(30,21).Value = "02/09/2023 14:51"
$ws.Cells.Item(30,22).Value = "https://www.betexplorer.com/football/denmark/2nd-division/thisted-fc-aarhus-fremad/C0kMLaoo/"
# Row 31 <- old row 32 data
$ws.Cells.Item(31,6).Value = "Thisted FC"
$ws.Cells.Item(31,7).Value = 0
$ws.Cells.Item(31,8).Value = "Skive"
$ws.Cells.Item(31,9).Value = 0
$ws.Cells.Item(31,10).Value = 2.16
$ws.Cells.Item(31,11).Value = "07/09/2023 07:12"
$ws.Cells.Item(31,12).Value = 2.41
$ws.Cells.Item(31,13).Value = "08/09/2023 18:19"
$ws.Cells.Item(31,14).Value = 3.37
$ws.Cells.Item(31,15).Value = "07/09/2023 07:12"
$ws.Cells.Item(31,16).Value = 3.39
$ws.Cells.Item(31,17).Value = "08/09/2023 18:56"
$ws.Cells.Item(31,18).Value = 2.9
$ws.Cells.Item(31,19).Value = "07/09/2023 07:12"
$ws.Cells.Item(31,20).Value = 2.77
$ws.Cells.Item(31,21).Value = "08/09/2023 18:56"
$ws.Cells.Item(31,22).Value = "https://www.betexplorer.com/football/denmark/2nd-division/thisted-fc-skive/MVWCrNhb/"
# Row 32 <- old row 33 data
$ws.Cells.Item(32,6).Value = "Roskilde"
$ws.Cells.Item(32,7).Value = 2
$ws.Cells.Item(32,8).Value = "AB Copenhagen"
$ws.Cells.Item(32,9).Value = 0
$ws.Cells.Item(32,10).Value = 2.15
$ws.Cells.Item(32,11).Value = "07/09/2023 07:12"
$ws.Cells.Item(32,12).Value = 2.25
$ws.Cells.Item(32,13).Value = "08/09/2023 18:58"
$ws.Cells.Item(32,14).Value = 3.42
$ws.Cells.Item(32,15).Value = "07/09/2023 07:12"
$ws.Cells.Item(32,16).Value = 3.34
$ws.Cells.Item(32,17).Value = "08/09/2023 18:58"
$ws.Cells.Item(32,18).Value = 2.79
$ws.Cells.Item(32,19).Value = "07/09/2023 07:12"
$ws.Cells.Item(32,20).Value = 3.05
$ws.Cells.Item(32,21).Value = "08/09/2023 18:58"
$ws.Cells.Item(32,22).Value = "https://www.betexplorer.com/football/denmark/2nd-division/roskilde-ab-copenhagen/v7s8q3wh/"
# Row 33 <- old row 31 data
$ws.Cells.Item(33,6).Value = "FA 2000"
$ws.Cells.Item(33,7).Value = 2
$ws.Cells.Item(33,8).Value = "Hellerup"
$ws.Cells.Item(33,9).Value = 0
$ws.Cells.Item(33,10).Value = 2.18
$ws.Cells.Item(33,11).Value = "08/09/2023 08:12"
$ws.Cells.Item(33,12).Value = 2.08
$ws.Cells.Item(33,13).Value = "08/09/2023 17:30"
$ws.Cells.Item(33,14).Value = 3.44
$ws.Cells.Item(33,15).Value = "08/09/2023 08:12"
$ws.Cells.Item(33,16).Value = 3.56
$ws.Cells.Item(33,17).Value = "08/09/2023 18:32"
$ws.Cells.Item(33,18).Value = 2.95
$ws.Cells.Item(33,19).Value = "08/09/2023 08:12"
$ws.Cells.Item(33,20).Value = 3.02
$ws.Cells.Item(33,21).Value = "08/09/2023 17:30"
$ws.Cells.Item(33,22).Value = "https://www.betexplorer.com/football/denmark/2nd-division/frederiksberg-alliancen-2000-hellerup/0Gr4pqOo/"

# --- Append new row 43 (Indice 42): Roskilde vs Hellerup, 22/09/2023 ---
# Copy formatting (bold/border style for the index cell, date-time number
# format for the match-date cell) from the previous last row (42) before
# writing the new values, so the appended row matches the sheet's existing
# per-column styling exactly.
$ws.Range("A42:V42").Copy()
$ws.Range("A43:V43").PasteSpecial(-4122)

$ws.Cells.Item(43,1).Value = 42
$ws.Cells.Item(43,2).Value = "denmark"
$ws.Cells.Item(43,3).Value = "2nd-division"
$ws.Cells.Item(43,4).Value = "2023-2024"
$ws.Cells.Item(43,5).Value = 45191.79166666666
$ws.Cells.Item(43,6).Value = "Roskilde"
$ws.Cells.Item(43,7).Value = 1
$ws.Cells.Item(43,8).Value = "Hellerup"
$ws.Cells.Item(43,9).Value = 1
$ws.Cells.Item(43,10).Value = 1.56
$ws.Cells.Item(43,11).Value = "21/09/2023 06:12"
$ws.Cells.Item(43,12).Value = 1.56
$ws.Cells.Item(43,13).Value = "22/09/2023 18:50"
$ws.Cells.Item(43,14).Value = 3.96
$ws.Cells.Item(43,15).Value = "21/09/2023 06:12"
$ws.Cells.Item(43,16).Value = 4.28
$ws.Cells.Item(43,17).Value = "22/09/2023 18:54"
$ws.Cells.Item(43,18).Value = 4.69
$ws.Cells.Item(43,19).Value = "21/09/2023 06:12"
$ws.Cells.Item(43,20).Value = 5.04
$ws.Cells.Item(43,21).Value = "22/09/2023 18:54"
$ws.Cells.Item(43,22).Value = "https://www.betexplorer.com/football/denmark/2nd-division/roskilde-hellerup/4WZ3thTR/"
